# Pooh Points: normal 20260207
# Applies the box-score refresh: status clock moved to 1:02 - 1st Half,
# several players' stat lines updated, "King Grace" split off Undrafted's
# row 16 into its own new row 17 (row 16 relabeled "Brandon Walker"), the
# old row 17 (Sergej Macura) shifting down to row 18, and the two
# corresponding OwnerTotals rollups bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# --- Update existing rows 2-16 cell-by-cell for changed values ---
$ws.Cells.Item(2,7).Value = "1:02 - 1st Half"

$ws.Cells.Item(3,7).Value = "1:02 - 1st Half"
$ws.Cells.Item(3,8).Value = 17
$ws.Cells.Item(3,9).Value = 13
$ws.Cells.Item(3,16).Value = 20
$ws.Cells.Item(3,21).Value = 3
$ws.Cells.Item(3,22).Value = 5

$ws.Cells.Item(4,7).Value = "1:02 - 1st Half"
$ws.Cells.Item(4,8).Value = 5
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 2
$ws.Cells.Item(4,16).Value = 20

$ws.Cells.Item(5,7).Value = "1:02 - 1st Half"
$ws.Cells.Item(5,16).Value = 13

$ws.Cells.Item(6,7).Value = "1:02 - 1st Half"

$ws.Cells.Item(7,7).Value = "1:02 - 1st Half"
$ws.Cells.Item(7,16).Value = 11

$ws.Cells.Item(8,7).Value = "1:02 - 1st Half"
$ws.Cells.Item(8,8).Value = 15
$ws.Cells.Item(8,9).Value = 12
$ws.Cells.Item(8,16).Value = 20
$ws.Cells.Item(8,17).Value = 5
$ws.Cells.Item(8,18).Value = 10
$ws.Cells.Item(8,19).Value = 1
$ws.Cells.Item(8,20).Value = 2

$ws.Cells.Item(9,7).Value = "1:02 - 1st Half"
$ws.Cells.Item(9,8).Value = 6
$ws.Cells.Item(9,14).Value = 3
$ws.Cells.Item(9,16).Value = 16

$ws.Cells.Item(10,7).Value = "1:02 - 1st Half"

$ws.Cells.Item(11,7).Value = "1:02 - 1st Half"
$ws.Cells.Item(11,8).Value = 12
$ws.Cells.Item(11,10).Value = 6
$ws.Cells.Item(11,16).Value = 16

$ws.Cells.Item(12,7).Value = "1:02 - 1st Half"
$ws.Cells.Item(12,16).Value = 11

$ws.Cells.Item(13,7).Value = "1:02 - 1st Half"

$ws.Cells.Item(14,7).Value = "1:02 - 1st Half"

$ws.Cells.Item(15,7).Value = "1:02 - 1st Half"

# Row 16 (Undrafted) becomes Brandon Walker with a trimmed stat line
$ws.Cells.Item(16,4).Value = "Brandon Walker"
$ws.Cells.Item(16,7).Value = "1:02 - 1st Half"
$ws.Cells.Item(16,10).Value = 0
$ws.Cells.Item(16,15).Value = 1
$ws.Cells.Item(16,16).Value = 1
$ws.Cells.Item(16,18).Value = 0
$ws.Cells.Item(16,20).Value = 0

# --- Insert a new row 17 for King Grace (pushes old row 17 -> row 18) ---
# Row 16 (still the original "King Grace" line at this point, before it is
# relabeled below) is copied down into the freshly-inserted row 17 first, so
# the new row inherits correct cell typing (notably the "date" column,
# which Excel would otherwise auto-parse into a date serial if assigned a
# plain "2026-02-07" string). Only the handful of cells that actually
# differ from row 16's original values are then overwritten.
$ws.Rows.Item(17).Insert()
$ws.Range("A16:V16").Copy($ws.Range("A17:V17"))
$ws.Cells.Item(17,7).Value = "1:02 - 1st Half"
$ws.Cells.Item(17,15).Value = 1
$ws.Cells.Item(17,16).Value = 6

# Row 18 is the old row 17 (Sergej Macura), shifted down by the insert;
# only its status clock needs refreshing.
$ws.Cells.Item(18,7).Value = "1:02 - 1st Half"

# --- Column G narrows from 18 to 17 characters wide ---
$ws.Columns.Item(7).ColumnWidth = 16.1

# --- OwnerTotals rollups for the two owners with starters so far ---
$ws2 = $wb.Worksheets.Item("OwnerTotals")
$ws2.Cells.Item(2,2).Value = 17
$ws2.Cells.Item(3,2).Value = 15
